$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data rows 8-11 (D:K input bits) - same pattern as existing row 7
$data = @(
    @(0,1,1,1,0,1,1,1),
    @(0,1,1,1,1,1,0,0),
    @(0,1,0,1,1,1,1,0),
    @(0,1,1,1,1,0,0,1)
)

$rowIdx = 8
foreach ($rowVals in $data) {
    $colIdx = 4 # D
    foreach ($v in $rowVals) {
        $ws.Cells.Item($rowIdx, $colIdx).Value = $v
        $colIdx++
    }
    $rowIdx++
}

# Column M: DEC2HEX formula filled down from row 7 through row 19
for ($r = 8; $r -le 19; $r++) {
    $ws.Range("M$r").Formula = "=DEC2HEX(K$r*1+J$r*2+I$r*4+H$r*8+G$r*16+F$r*32+E$r*64+D$r*128)"
}

# Columns O,P,Q,R,S,T,U,V,X for rows 8-11 mirror row 7's formulas
for ($r = 8; $r -le 11; $r++) {
    $ws.Range("O$r").Formula = "=G$r"
    $ws.Range("P$r").Formula = "=H$r"
    $ws.Range("Q$r").Formula = "=I$r"
    $ws.Range("R$r").Formula = "=J$r"
    $ws.Range("S$r").Formula = "=K$r"
    $ws.Range("T$r").Formula = "=D$r"
    $ws.Range("U$r").Formula = "=F$r"
    $ws.Range("V$r").Formula = "=E$r"
    $ws.Range("X$r").Formula = "=DEC2HEX(V$r*1+U$r*2+T$r*4+S$r*8+R$r*16+Q$r*32+P$r*64+O$r*128)"
}

# Update selection to match target state
$ws.Range("M11:W11").Select() | Out-Null

$wb.Save()
